{"js": "// Load all paragraphs in the document body so we can inspect their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) Title: \"LOB1211 -  Polui\u00e7\u00e3o Ambiental II\" -> \"LOB1211 -  Polui\u00e7\u00e3o Atmosf\u00e9rica\"\n// 2) \"Curso (semestre ideal): EA (6)\" -> \"Curso (semestre ideal): EA (7)\"\n// Use a scoped search + insertText(\"Replace\") on just the matching run text so\n// surrounding runs/marks (e.g. the <w:br/> siblings on the credits paragraph)\n// are left untouched.\nfor (const p of items) {\n  if (p.text.indexOf(\"LOB1211 -  Polui\u00e7\u00e3o Ambiental II\") !== -1) {\n    const found = p.search(\"LOB1211 -  Polui\u00e7\u00e3o Ambiental II\", { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(\"LOB1211 -  Polui\u00e7\u00e3o Atmosf\u00e9rica\", \"Replace\");\n    }\n  }\n  if (p.text.indexOf(\"Curso (semestre ideal): EA (6)\") !== -1) {\n    const found = p.search(\"Curso (semestre ideal): EA (6)\", { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(\"Curso (semestre ideal): EA (7)\", \"Replace\");\n    }\n  }\n}\nawait context.sync();\n\n// 3) Remove the \"Requisitos\" Heading2 paragraph and the ListBullet paragraph\n// that follows it (the three \"(Requisito)\" lines), leaving the Bibliografia\n// paragraph directly followed by the section end.\nconst toDelete = [];\nfor (const p of items) {\n  const t = p.text;\n  if (t === \"Requisitos\" || t.indexOf(\"LOB1207 -  Polui\u00e7\u00e3o Ambiental I\") !== -1) {\n    toDelete.push(p);\n  }\n}\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title: \"LOB1211 -  Polui\u00e7\u00e3o Ambiental II\" -> \"LOB1211 -  Polui\u00e7\u00e3o Atmosf\u00e9rica\"\n$find1 = $d.Content.Find\n$find1.Text = \"LOB1211 -  Polui\u00e7\u00e3o Ambiental II\"\n$find1.Replacement.Text = \"LOB1211 -  Polui\u00e7\u00e3o Atmosf\u00e9rica\"\n$find1.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) \"Curso (semestre ideal): EA (6)\" -> \"Curso (semestre ideal): EA (7)\"\n$find2 = $d.Content.Find\n$find2.Text = \"Curso (semestre ideal): EA (6)\"\n$find2.Replacement.Text = \"Curso (semestre ideal): EA (7)\"\n$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3) Remove the \"Requisitos\" Heading2 paragraph and the ListBullet paragraph\n# that follows it (the three \"(Requisito)\" lines). These are the last two\n# paragraphs in the body (right before the section properties), so deleting\n# the range from the start of \"Requisitos\" through the end of the document's\n# last paragraph removes them both and leaves Bibliografia directly followed\n# by the section end.\n$startPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Requisitos\") {\n        $startPara = $p\n        break\n    }\n}\n\nif ($startPara -ne $null) {\n    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $rng = $d.Range($startPara.Range.Start, $lastPara.Range.End)\n    $rng.Delete()\n}\n"}
